# Insert a new daily price-report row for Mango at row 445 (shifting the
# existing rows 445-495 down to 446-496), then populate the new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(445).Insert()

$ws.Range("A445").Value = 10
$ws.Range("B445").Value = "Vega Modelo de Temuco"
$ws.Range("C445").Value = "La Araucanía"
$ws.Range("D445").Value = 44946
$ws.Range("E445").Value = 9
$ws.Range("F445").Value = "Fruta"
$ws.Range("G445").Value = 100108
$ws.Range("H445").Value = "Tropicales y subtropicales"
$ws.Range("I445").Value = 100108002
$ws.Range("J445").Value = "Mango"
$ws.Range("K445").Value = "Sin especificar"
$ws.Range("L445").Value = "Primera"
$ws.Range("M445").Value = 195
$ws.Range("N445").Value = 7500
$ws.Range("O445").Value = 7500
$ws.Range("P445").Value = 7500
$ws.Range("Q445").Value = "$/bandeja 4 kilos"
$ws.Range("R445").Value = "Perú"
$ws.Range("S445").Value = 1875
$ws.Range("T445").Value = 4
